# Auto-generated edit script: updates cryptos list per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-ExactText {
    param($cellRef, $text)
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.ClearFormats()
}

Set-ExactText 'D2' '63.106.31'
Set-ExactText 'E2' '  -2.04%  '
Set-ExactText 'D3' '2.624.11'
Set-ExactText 'E3' '  -1.97%  '
Set-ExactText 'D4' '0.999'
Set-ExactText 'E4' '  -0.09%  '
Set-ExactText 'D5' '604.63'
Set-ExactText 'E5' '  +1.38%  '
Set-ExactText 'D6' '146.29'
Set-ExactText 'E6' '  -1.13%  '
Set-ExactText 'D7' '0.999'
Set-ExactText 'E7' '  -0.08%  '
Set-ExactText 'E8' '  -1.32%  '
Set-ExactText 'D9' '2.623.71'
Set-ExactText 'E9' '  -1.94%  '
Set-ExactText 'E10' '  -0.42%  '
Set-ExactText 'E11' '  -1.07%  '
Set-ExactText 'E12' '  -0.05%  '
Set-ExactText 'E13' '  +0.87%  '
Set-ExactText 'D14' '27.18'
Set-ExactText 'E14' '  -2.84%  '
Set-ExactText 'D15' '3.091.10'
Set-ExactText 'E15' '  -1.94%  '
Set-ExactText 'D16' '62.900.94'
Set-ExactText 'E16' '  -2.21%  '
Set-ExactText 'D17' '0.0000145'
Set-ExactText 'E17' '  -2.17%  '
Set-ExactText 'D18' '2.629.88'
Set-ExactText 'E18' '  -2.06%  '
Set-ExactText 'D19' '11.28'
Set-ExactText 'E19' '  -1.35%  '
Set-ExactText 'D20' '4.49'
Set-ExactText 'E20' '  +1.63%  '
Set-ExactText 'D21' '339.99'
Set-ExactText 'E21' '  -2.00%  '
Set-ExactText 'E22' '  -0.64%  '
Set-ExactText 'E23' '  -0.08%  '
Set-ExactText 'E24' '  -4.63%  '
Set-ExactText 'D25' '66.47'
Set-ExactText 'E25' '  -3.44%  '
Set-ExactText 'D26' '1.62'
Set-ExactText 'E26' '  -3.12%  '
Set-ExactText 'E27' '  -5.04%  '
Set-ExactText 'E28' '  +1.26%  '
Set-ExactText 'E29' '  -2.71%  '
Set-ExactText 'B30' 'Bittensor'
Set-ExactText 'C30' 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-ExactText 'D30' '541.50'
Set-ExactText 'E30' '  +2.28%  '
Set-ExactText 'B31' 'Binance-PegBSC-USD'
Set-ExactText 'C31' 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-ExactText 'D31' '1.00'
Set-ExactText 'E31' '  +0.11%  '
Set-ExactText 'B32' 'Aptos'
Set-ExactText 'C32' 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-ExactText 'D32' '7.91'
Set-ExactText 'E32' '  -1.40%  '
Set-ExactText 'E33' '  +1.23%  '
Set-ExactText 'D34' '1.75'
Set-ExactText 'E34' '  -2.40%  '
Set-ExactText 'D35' '0.0₃0803'
Set-ExactText 'E35' '  -2.99%  '
Set-ExactText 'E36' '  +10.80%  '
Set-ExactText 'D37' '169.19'
Set-ExactText 'E37' '  -3.62%  '
Set-ExactText 'E38' '  -0.18%  '
Set-ExactText 'D39' '0.403'
Set-ExactText 'E39' '  -0.30%  '
Set-ExactText 'E40' '  -1.74%  '
Set-ExactText 'E41' '  +5.57%  '
Set-ExactText 'E42' '  +0.00%  '
Set-ExactText 'D43' '169.53'
Set-ExactText 'E43' '  -2.04%  '
Set-ExactText 'E44' '  -1.18%  '
Set-ExactText 'D45' '22.30'
Set-ExactText 'E45' '  +1.97%  '
Set-ExactText 'D46' '0.0566'
Set-ExactText 'E46' '  +2.67%  '
Set-ExactText 'E47' '  -1.89%  '
Set-ExactText 'E48' '  -0.86%  '
Set-ExactText 'E49' '  -0.43%  '
Set-ExactText 'D50' '18.45'
Set-ExactText 'E50' '  -2.27%  '
Set-ExactText 'E51' '  +0.09%  '

Write-Host "Applied $($wb.Worksheets.Count) sheet(s); updated 84 cells."
